$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text edits in header strings ---
$ws.Range("A8").Replace("17", "18") | Out-Null
$ws.Range("C9").Replace("4/24/2023", "5/1/2023") | Out-Null
$ws.Range("C9").Replace("4/30/2023", "5/7/2023") | Out-Null

# --- Reference cells for formats (unchanged cells used as PasteSpecial format sources) ---
$refText = $ws.Range("C14")   # style 14 (General/text placeholder)
$refInt  = $ws.Range("I14")   # style 15 (#,##0 integer)
$refPct  = $ws.Range("K14")   # style 16 (#,##0.0 percent-like)

# --- Data cell edits, rows 14-30 ---
# Row 14
$ws.Range("N14").Value = -84.615384615384

# Row 15
$ws.Range("C15").Value = "'0"
$refText.Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("D15").Value = "'0"
$refText.Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("E15").Value = "***.*"
$refText.Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -40
$ws.Range("N15").Value = 100

# Row 16
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -14.285714285714
$ws.Range("F16").Value = 29
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = 38.095238095238
$ws.Range("I16").Value = 126
$ws.Range("J16").Value = 103
$ws.Range("K16").Value = 22.330097087378
$ws.Range("L16").Value = 85.294117647058
$ws.Range("M16").Value = -4.545454545454
$ws.Range("N16").Value = -77.777777777777

# Row 17
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -30
$ws.Range("F17").Value = 46
$ws.Range("G17").Value = 42
$ws.Range("H17").Value = 9.523809523809
$ws.Range("I17").Value = 209
$ws.Range("J17").Value = 164
$ws.Range("K17").Value = 27.439024390243
$ws.Range("L17").Value = 65.873015873015
$ws.Range("M17").Value = 120
$ws.Range("N17").Value = 43.150684931506

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = -23.529411764705
$ws.Range("I18").Value = 74
$ws.Range("J18").Value = 57
$ws.Range("K18").Value = 29.824561403508
$ws.Range("L18").Value = 19.354838709677
$ws.Range("M18").Value = -45.185185185185
$ws.Range("N18").Value = -91.158900836320

# Row 19
$ws.Range("F19").Value = 72
$ws.Range("G19").Value = 72
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 303
$ws.Range("J19").Value = 446
$ws.Range("K19").Value = -32.062780269058
$ws.Range("L19").Value = 103.355704697987
$ws.Range("M19").Value = 56.994818652849
$ws.Range("N19").Value = -20.472440944881

# Row 20
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -20
$ws.Range("F20").Value = 20
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = 11.111111111111
$ws.Range("I20").Value = 95
$ws.Range("J20").Value = 64
$ws.Range("K20").Value = 48.4375
$ws.Range("L20").Value = 102.127659574468
$ws.Range("M20").Value = 72.727272727272
$ws.Range("N20").Value = -86.787204450625

# Row 21
$ws.Range("C21").Value = 35
$ws.Range("D21").Value = 42
$ws.Range("E21").Value = -16.666666666666
$ws.Range("F21").Value = 183
$ws.Range("G21").Value = 175
$ws.Range("H21").Value = 4.571428571428
$ws.Range("I21").Value = 819
$ws.Range("J21").Value = 847
$ws.Range("K21").Value = -3.305785123966
$ws.Range("L21").Value = 75
$ws.Range("M21").Value = 31.884057971014
$ws.Range("N21").Value = -69.302848575712

# Row 22
$ws.Range("D22").Value = "'0"
$refText.Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").Value = "***.*"
$refText.Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("F22").Value = "'0"
$refText.Copy() | Out-Null
$ws.Range("F22").PasteSpecial(-4122) | Out-Null
$ws.Range("H22").Value = -100
$ws.Range("M22").Value = 7.692307692307

# Row 24
$ws.Range("C24").Value = 48
$ws.Range("D24").Value = 43
$ws.Range("E24").Value = 11.627906976744
$ws.Range("F24").Value = 207
$ws.Range("G24").Value = 167
$ws.Range("H24").Value = 23.952095808383
$ws.Range("I24").Value = 971
$ws.Range("J24").Value = 805
$ws.Range("K24").Value = 20.621118012422
$ws.Range("L24").Value = 77.189781021897
$ws.Range("M24").Value = 86.015325670498

# Row 25
$ws.Range("C25").Value = 28
$ws.Range("D25").Value = 20
$ws.Range("E25").Value = 40
$ws.Range("F25").Value = 94
$ws.Range("G25").Value = 59
$ws.Range("H25").Value = 59.322033898305
$ws.Range("I25").Value = 370
$ws.Range("J25").Value = 261
$ws.Range("K25").Value = 41.762452107279
$ws.Range("L25").Value = 58.798283261802
$ws.Range("M25").Value = 77.033492822966

# Row 26
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = "'0"
$refText.Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null
$ws.Range("E26").Value = "***.*"
$refText.Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4122) | Out-Null
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 16
$ws.Range("K26").Value = -5.882352941176
$ws.Range("L26").Value = -20

# Row 27
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 1
$refInt.Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("E27").Value = 100
$refPct.Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("F27").Value = 6
$ws.Range("H27").Value = 20
$ws.Range("I27").Value = 44
$ws.Range("J27").Value = 33
$ws.Range("K27").Value = 33.333333333333
$ws.Range("L27").Value = 22.222222222222

# Row 28
$ws.Range("N28").Value = -88

# Row 29
$ws.Range("N29").Value = -91.666666666666

# Row 30
$ws.Range("D30").Value = 1
$refInt.Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null
$ws.Range("E30").Value = -100
$refPct.Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4122) | Out-Null
$ws.Range("G30").Value = 1
$refInt.Copy() | Out-Null
$ws.Range("G30").PasteSpecial(-4122) | Out-Null
$ws.Range("H30").Value = -100
$refPct.Copy() | Out-Null
$ws.Range("H30").PasteSpecial(-4122) | Out-Null
$ws.Range("J30").Value = 4
$ws.Range("K30").Value = 0

$excel.CutCopyMode = $false